$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value that was bulk-updated for every
# data row (rows 2-480) from 45172 (2023-09-03) to 45175 (2023-09-06).
$ws.Range("C2:C480").Value = 45175
